$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.931.74'
$ws.Range("E2").Value = '  -1.28%  '
$ws.Range("D3").Value = '1.637.40'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '215.37'
$ws.Range("E6").Value = '  -0.34%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  -0.88%  '
$ws.Range("E9").Value = '  -0.10%  '
$ws.Range("D10").Value = '19.60'
$ws.Range("E10").Value = '  -2.26%  '
$ws.Range("D11").Value = '0.0795'
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").Value = '1.864.58'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").Value = '1.636.51'
$ws.Range("E14").Value = '  -0.06%  '
$ws.Range("E15").Value = '  -1.16%  '
$ws.Range("E16").Value = '  -0.40%  '
$ws.Range("D17").Value = '62.96'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").Value = '25.895.67'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("E19").Value = '  +0.22%  '
$ws.Range("D20").Value = '192.79'
$ws.Range("E20").Value = '  -1.66%  '
$ws.Range("E21").Value = '  -2.11%  '
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("E24").Value = '  +0.95%  '
$ws.Range("D25").Value = '144.16'
$ws.Range("E25").Value = '  +0.43%  '
$ws.Range("E27").Value = '  +3.73%  '
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("D29").Value = '15.55'
$ws.Range("E29").Value = '  -0.51%  '
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("D31").Value = '0.0502'
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = '3.30'
$ws.Range("E32").Value = '  -1.45%  '
$ws.Range("E33").Value = '  -0.58%  '
$ws.Range("D34").Value = '1.53'
$ws.Range("E34").Value = '  -4.69%  '
$ws.Range("D35").Value = '2.44'
$ws.Range("E35").Value = '  +1.62%  '
$ws.Range("D36").Value = '0.901'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").Value = '1.136.09'
$ws.Range("E37").Value = '  -0.30%  '
$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  -1.51%  '
$ws.Range("E40").Value = '  +0.19%  '
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").Value = '99.19'
$ws.Range("E42").Value = '  -1.13%  '
$ws.Range("D43").Value = '0.797'
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("D44").Value = '1.774.45'
$ws.Range("E44").Value = '  -0.48%  '
$ws.Range("D45").Value = '0.0₆0114'
$ws.Range("E45").Value = '  +1.59%  '
$ws.Range("D46").Value = '56.64'
$ws.Range("E46").Value = '  +0.52%  '
$ws.Range("D47").Value = '0.0530'
$ws.Range("E47").Value = '  +2.70%  '
$ws.Range("E48").Value = '  +0.31%  '
$ws.Range("D49").Value = '7.68'
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("E50").Value = '  -1.01%  '
$ws.Range("D51").Value = '0.0965'
$ws.Range("E51").Value = '  -0.91%  '
